$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, pushing the existing "rent" row down to row 3
$ws.Rows.Item(2).Insert()

# New row 2: dnace / 300 / 2026-02-03 05:30:10 (serial 46056.22928240741)
$ws.Cells.Item(2, 1).Value = "dnace"
$ws.Cells.Item(2, 2).Value = 300
$ws.Cells.Item(2, 3).Value = 46056.22928240741

# New row 4: Freelance Project / 56200 / 2026-01-28 23:30:10 (serial 46050.97928240741)
$ws.Cells.Item(4, 1).Value = "Freelance Project"
$ws.Cells.Item(4, 2).Value = 56200
$ws.Cells.Item(4, 3).Value = 46050.97928240741

# Copy the date cell format from row 3 (original style) onto the two new date cells
# so they reuse the same style index instead of minting new ones.
$ws.Cells.Item(3, 3).Copy()
$ws.Cells.Item(2, 3).PasteSpecial(-4122)
$ws.Cells.Item(4, 3).PasteSpecial(-4122)
$excel.CutCopyMode = $false
